# Add two new visitor records (rows 3 and 4) to the "Visitors" worksheet,
# matching the rows already present (Name, Phone, Reason, Status, Date, Photo URL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "dcsdc"
$ws.Range("B3").Value = "dfvdfv"
$ws.Range("C3").Value = "fvdfa"
$ws.Range("D3").Value = "approved"
$ws.Range("E3").Value = "2/7/2025, 9:21:18 pm"
$ws.Range("F3").Value = "https://res.cloudinary.com/drdw2abup/image/upload/v1751471478/visitors/iexcozqvsik1snensf9n.jpg"

# Row 4
$ws.Range("A4").Value = "fefaef"
$ws.Range("B4").Value = "gsergrse"
$ws.Range("C4").Value = "gtsrtg"
$ws.Range("D4").Value = "approved"
$ws.Range("E4").Value = "2/7/2025, 9:28:47 pm"
$ws.Range("F4").Value = "https://res.cloudinary.com/drdw2abup/image/upload/v1751471927/visitors/q1cx4nvgvwz4asjemjf1.jpg"
